# Added udp sender run 1 to 5
$wb = $excel.ActiveWorkbook

# --- Rename the second sheet to reflect the new UDP Sender example ---
$wsJson = $wb.Worksheets.Item(1)
$wsUdp  = $wb.Worksheets.Item(2)
$wsUdp.Name = "UDP Sender Code Beispiel"

# --- Enter the five new UDP sender measurement runs (rows 2-6) ---
# Column B = Index values, column C = CompilerErr values
$wsUdp.Range("B2").Value = 60
$wsUdp.Range("C2").Value = 0

$wsUdp.Range("B3").Value = 59
$wsUdp.Range("C3").Value = 0

$wsUdp.Range("B4").Value = 60
$wsUdp.Range("C4").Value = 0

$wsUdp.Range("B5").Value = 65
$wsUdp.Range("C5").Value = 6

$wsUdp.Range("B6").Value = 63
$wsUdp.Range("C6").Value = 0

# --- Restore the selections left on each sheet ---
$wsJson.Activate() | Out-Null
$wsJson.Range("F34").Select() | Out-Null

$wsUdp.Activate() | Out-Null
$wsUdp.Range("F6").Select() | Out-Null
